$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 - California
$ws.Range("B2").Value = 45
$ws.Range("C2").Value = 9
$ws.Range("D2").Value = 16
$ws.Range("E2").Value = 7

# Row 3 - Los Angeles
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 11
$ws.Range("E3").Value = 3

# Row 4 - San Diego
$ws.Range("B4").Value = 50
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = 6
$ws.Range("E4").Value = 16

# Update the selected cell to reflect the saved view state
$ws.Range("F3").Select()
